# Daily attendance processing - 2025-11-07 15:23:02
# Normalizes the "Recorded By" (column G) audit-trail strings on the
# "Session Analysis Results" sheet: for every row whose recorded-by value
# lists two or more comma-separated contributors, the first two entries
# are swapped (any additional entries keep their original order/position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ($text -ne $null -and $text -ne "") {
        $parts = $text -split ', '
        if ($parts.Length -ge 2) {
            $first = $parts[0]
            $second = $parts[1]
            $parts[0] = $second
            $parts[1] = $first
            $cell.Value = ($parts -join ', ')
        }
    }
}
